$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 538.44446
$ws.Range("I9").Value = 468.69232
$ws.Range("J9").Value = 719.8
$ws.Range("K9").Value = 468.69232
$ws.Range("L9").Value = 719.8
$ws.Range("M9").Value = -299.69232
$ws.Range("N9").Value = -1057.8

$ws.Range("H32").Value = 8198
$ws.Range("I32").Value = 12263.286
$ws.Range("J32").Value = 6419.4375
$ws.Range("K32").Value = 12263.286
$ws.Range("L32").Value = 6419.4375
$ws.Range("M32").Value = -11937.286
$ws.Range("N32").Value = -7071.4375

$ws.Range("H33").Value = 537.25
$ws.Range("J33").Value = 630
$ws.Range("L33").Value = 630
$ws.Range("N33").Value = -1088

$ws.Range("H40").Value = 2883.8
$ws.Range("J40").Value = 3473
$ws.Range("L40").Value = 3473
$ws.Range("N40").Value = -3823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41468.27
$ws.Range("I32").Value = 44757.332
$ws.Range("K32").Value = 44757.332
$ws.Range("M32").Value = -44470.332

$ws.Range("H74").Value = 218812.72
$ws.Range("I74").Value = 286269.72
$ws.Range("J74").Value = 16441.715
$ws.Range("K74").Value = 286269.72
$ws.Range("L74").Value = 16441.715
$ws.Range("M74").Value = -285395.72
$ws.Range("N74").Value = -18189.715

$ws.Range("H77").Value = 218812.72
$ws.Range("I77").Value = 286269.72
$ws.Range("J77").Value = 16441.715
$ws.Range("K77").Value = 1431348.6
$ws.Range("L77").Value = 82208.575
$ws.Range("M77").Value = -1426980.6
$ws.Range("N77").Value = -90944.575

$ws.Range("H122").Value = 3455.5334
$ws.Range("I122").Value = 3342
$ws.Range("J122").Value = 3828.5715
$ws.Range("K122").Value = 10026
$ws.Range("L122").Value = 11485.7145
$ws.Range("M122").Value = -7576
$ws.Range("N122").Value = -16385.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19369.111
$ws.Range("I20").Value = 33056.3
$ws.Range("J20").Value = 2260.125
$ws.Range("K20").Value = 33056.3
$ws.Range("L20").Value = 2260.125
$ws.Range("M20").Value = -32809.3
$ws.Range("N20").Value = -2754.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 950.2727
$ws.Range("I5").Value = 759.64703
$ws.Range("J5").Value = 1598.4
$ws.Range("K5").Value = 2278.94109
$ws.Range("L5").Value = 4795.200000000001
$ws.Range("M5").Value = -2166.94109
$ws.Range("N5").Value = -5019.200000000001

$ws.Range("H23").Value = 366.7143
$ws.Range("J23").Value = 592.7143
$ws.Range("L23").Value = 1778.1429
$ws.Range("N23").Value = -2248.1429

$ws.Range("H26").Value = 61.666668
$ws.Range("I26").Value = 48.125
$ws.Range("K26").Value = 144.375
$ws.Range("M26").Value = 143.625

$ws.Range("H34").Value = 7749.8335
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2916

$ws.Range("H40").Value = 345.33334
$ws.Range("I40").Value = 68
$ws.Range("J40").Value = 900
$ws.Range("K40").Value = 272
$ws.Range("L40").Value = 3600
$ws.Range("M40").Value = -203
$ws.Range("N40").Value = -3738

$ws.Range("H45").Value = 2064.125
$ws.Range("I45").Value = 1318.8
$ws.Range("J45").Value = 3306.3333
$ws.Range("K45").Value = 3956.4
$ws.Range("L45").Value = 9918.999899999999
$ws.Range("M45").Value = -3424.4
$ws.Range("N45").Value = -10982.9999

$ws.Range("H59").Value = 4032.4075
$ws.Range("I59").Value = 937.5
$ws.Range("J59").Value = 4280
$ws.Range("K59").Value = 2812.5
$ws.Range("L59").Value = 12840
$ws.Range("M59").Value = -2272.5
$ws.Range("N59").Value = -13920

$ws.Range("H131").Value = 2767.8572
$ws.Range("I131").Value = 3002.2666
$ws.Range("J131").Value = 2181.8333
$ws.Range("K131").Value = 9006.799800000001
$ws.Range("L131").Value = 6545.499899999999
$ws.Range("M131").Value = -3966.799800000001
$ws.Range("N131").Value = -16625.4999

$ws.Range("H133").Value = 8429
$ws.Range("J133").Value = 8429
$ws.Range("L133").Value = 25287
$ws.Range("N133").Value = -35407

$ws.Range("H135").Value = 950.2727
$ws.Range("I135").Value = 759.64703
$ws.Range("J135").Value = 1598.4
$ws.Range("K135").Value = 6836.82327
$ws.Range("L135").Value = 14385.6
$ws.Range("M135").Value = -4301.82327
$ws.Range("N135").Value = -19455.6

$ws.Range("H139").Value = 13771.8
$ws.Range("I139").Value = 13771.8
$ws.Range("K139").Value = 41315.39999999999
$ws.Range("M139").Value = -36175.39999999999

$ws.Range("H140").Value = 7255
$ws.Range("I140").Value = 7255
$ws.Range("K140").Value = 21765
$ws.Range("M140").Value = -16585

$ws.Range("H141").Value = 6264.3335
$ws.Range("I141").Value = 5625.5713
$ws.Range("J141").Value = 8500
$ws.Range("K141").Value = 16876.7139
$ws.Range("L141").Value = 25500
$ws.Range("M141").Value = -11696.7139
$ws.Range("N141").Value = -35860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H70").Value = 6244.1665
$ws.Range("I70").Value = 6186.2856
$ws.Range("J70").Value = 6325.2
$ws.Range("K70").Value = 6186.2856
$ws.Range("L70").Value = 6325.2
$ws.Range("M70").Value = -5916.2856
$ws.Range("N70").Value = -6865.2

$ws.Range("H73").Value = 6244.1665
$ws.Range("I73").Value = 6186.2856
$ws.Range("J73").Value = 6325.2
$ws.Range("K73").Value = 6186.2856
$ws.Range("L73").Value = 6325.2
$ws.Range("M73").Value = -5250.2856
$ws.Range("N73").Value = -8197.200000000001

$ws.Range("H80").Value = 19457
$ws.Range("I80").Value = 13999.667
$ws.Range("J80").Value = 23550
$ws.Range("K80").Value = 13999.667
$ws.Range("L80").Value = 23550
$ws.Range("M80").Value = -13001.667
$ws.Range("N80").Value = -25546

$ws.Range("H83").Value = 19457
$ws.Range("I83").Value = 13999.667
$ws.Range("J83").Value = 23550
$ws.Range("K83").Value = 69998.33499999999
$ws.Range("L83").Value = 117750
$ws.Range("M83").Value = -65006.33499999999
$ws.Range("N83").Value = -127734

$ws.Range("H122").Value = 2492.3823
$ws.Range("I122").Value = 2288.8076
$ws.Range("J122").Value = 3154
$ws.Range("K122").Value = 6866.4228
$ws.Range("L122").Value = 9462
$ws.Range("M122").Value = -4416.4228
$ws.Range("N122").Value = -14362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4166.6665
$ws.Range("I122").Value = 4166.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12499.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10049.9995
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 13374.25
$ws.Range("J49").Value = 13374.25
$ws.Range("L49").Value = 13374.25
$ws.Range("N49").Value = -13834.25

$ws.Range("H81").Value = 10808.5
$ws.Range("J81").Value = 5097.5
$ws.Range("L81").Value = 10195
$ws.Range("N81").Value = -12317

$ws.Range("H84").Value = 10808.5
$ws.Range("J84").Value = 5097.5
$ws.Range("L84").Value = 50975
$ws.Range("N84").Value = -61583

$ws.Range("H136").Value = 33732.633
$ws.Range("I136").Value = 47110.152
$ws.Range("K136").Value = 141330.456
$ws.Range("M136").Value = -138780.456
